$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update quantity for the LED row (row 4) from 2 to 1.
# Dependent formulas (F4 = D4*E4, F6 = SUM(F2:F4)) recalc automatically.
$ws.Range("D4").Value = 1

# Recalculate to ensure formula results are refreshed.
$excel.Calculate()

# Update the active selection to match the authored change.
$ws.Activate()
$ws.Range("D5").Select()
